$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update the generation "Date" property (row 8, column B) ---
# This is a plain text timestamp string, so a normal value assignment is fine.
$ws.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# --- Fill in the "Case Sensitive" property value (row 15, column B) ---
# The target value is the literal text "true" (not the Boolean TRUE), so a
# direct $cell.Value = "true" assignment can't be used here (Excel's input
# parser would coerce a bare "true"/"false" string into a Boolean). To force
# plain text while keeping the cell's existing style untouched, stage the
# text in a scratch cell using a leading apostrophe (the normal Excel way of
# entering literal text that looks like another type), copy it across with
# PasteSpecial (values only, so formatting on B15 is left alone), then wipe
# the scratch cell.
$scratch = $ws.Range("Z1")
$scratch.Value = "'true"
$scratch.Copy()
$ws.Range("B15").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
